$d = $word.ActiveDocument

# Locate the paragraph that holds the "USCOTS 2013" run (right-aligned,
# 72pt grey heading line) and the paragraph right after it (the blank
# 96pt line that follows). We find them by scanning for the "USCOTS 2013"
# text rather than hard-coding paragraph indices, so the script is robust
# to the surrounding content.
$uscotsIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "USCOTS 2013") {
        $uscotsIndex = $i
        break
    }
}

if ($uscotsIndex -gt 0) {
    $uscotsPara = $d.Paragraphs($uscotsIndex)
    $uscotsRange = $uscotsPara.Range

    # Remove the "USCOTS 2013" run entirely (delete up to, but excluding,
    # the paragraph mark so the empty paragraph itself is preserved).
    $deleteRange = $d.Range($uscotsRange.Start, $uscotsRange.End - 1)
    $deleteRange.Delete()

    # The next paragraph (originally empty) is where the "_GoBack" bookmark
    # should now live, marking the most recent edit location.
    $nextPara = $d.Paragraphs($uscotsIndex + 1)
    $nextRange = $nextPara.Range
    $nextRange.Collapse(1)   # wdCollapseStart
    $d.Bookmarks.Add("_GoBack", $nextRange) | Out-Null
}
